{"js": "// Bump stack-trace line numbers in the stored exception text (M2Doc 3.0.0 -> 3.1.0)\n// and add the extra \"RunBefores.evaluate\" frame that appears with the new version,\n// matching the unified diff.\n\nasync function replaceAll(context, searchText, replaceText) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Simple line-number bumps (each search string is unique in the document).\nawait replaceAll(context, \"M2DocEvaluator.caseQuery(M2DocEvaluator.java:559)\", \"M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)\");\n\n// This one repeats 3 times in the stack trace; replace every occurrence.\nawait replaceAll(context, \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)\", \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)\");\n\nawait replaceAll(context, \"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)\", \"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)\");\nawait replaceAll(context, \"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)\", \"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)\");\nawait replaceAll(context, \"M2DocEvaluator.generate(M2DocEvaluator.java:276)\", \"M2DocEvaluator.generate(M2DocEvaluator.java:281)\");\nawait replaceAll(context, \"M2DocUtils.generate(M2DocUtils.java:694)\", \"M2DocUtils.generate(M2DocUtils.java:805)\");\nawait replaceAll(\n  context,\n  \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\",\n  \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)\"\n);\nawait replaceAll(\n  context,\n  \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\",\n  \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)\"\n);\n\n// New JUnit frame inserted right before the (unique) RunAfters line that\n// immediately follows the \"ParentRunner$2.evaluate(ParentRunner.java:268)\" frame.\nawait replaceAll(\n  context,\n  \"ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\",\n  \"ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\\n\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\"\n);\n", "ps1": "# Bump stack-trace line numbers in the stored exception text (M2Doc 3.0.0 -> 3.1.0)\n# and add the extra \"RunBefores.evaluate\" frame that appears with the new version,\n# matching the unified diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($find, [string]$searchText, [string]$replaceText) {\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n$find = $d.Content.Find\n\n# Simple line-number bumps (each search string is unique in the document).\nReplace-All $find \"M2DocEvaluator.caseQuery(M2DocEvaluator.java:559)\" \"M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)\"\n\n# This one repeats 3 times in the stack trace; wdReplaceAll handles every occurrence.\nReplace-All $find \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)\" \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)\"\n\nReplace-All $find \"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)\" \"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)\"\nReplace-All $find \"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)\" \"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)\"\nReplace-All $find \"M2DocEvaluator.generate(M2DocEvaluator.java:276)\" \"M2DocEvaluator.generate(M2DocEvaluator.java:281)\"\nReplace-All $find \"M2DocUtils.generate(M2DocUtils.java:694)\" \"M2DocUtils.generate(M2DocUtils.java:805)\"\nReplace-All $find \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\" \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)\"\nReplace-All $find \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\" \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)\"\n\n# New JUnit frame inserted right before the (unique) RunAfters line that\n# immediately follows the \"ParentRunner$2.evaluate(ParentRunner.java:268)\" frame.\n$searchText = \"ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\"\n$replaceText = \"ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\"\nReplace-All $find $searchText $replaceText\n"}
